# Assignment2/Benchmarks.xlsx - "Fix pre-allocation bug, rerun benchmarks,
# modify probability assignment approach"
#
# The only user-visible/content changes in this revision are on Sheet1:
#   - The note in cell B13 ("Computed using 3.1 GHz Intel Core i7 quad-core
#     processor (4 physical, 8 logical)") gains a trailing space, which in
#     turn causes Excel to re-emit the shared-string table with that entry
#     swapped to the end (the "Test 3 (read-dominated)" label in F1 keeps
#     its displayed text but now points at the other shared-string index).
#   - The worksheet's last selected cell moves from B14 to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-enter the footnote text with an added trailing space.
$ws.Range("B13").Value = "Computed using 3.1 GHz Intel Core i7 quad-core processor (4 physical, 8 logical) "

# Leave the selection on B13 (matches the saved cursor position in the diff).
$ws.Range("B13").Select() | Out-Null
